$d = $word.ActiveDocument

# 1. Paragraph "Кнопка: «Отправить»" - highlight whole paragraph (incl. paragraph mark) darkCyan
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Кнопка:*" -and $t -like "*«Отправить»*") {
        $p.Range.Font.HighlightColorIndex = "darkCyan"
    }
}

# 2. Paragraph "5. Вы блестяще справились с Квестом!..." - highlight only the
#    sentence (not the leading "5. " nor the trailing space) darkCyan
$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("Вы блестяще справились с Квестом! Мы рады видеть Вас в качестве клиента и перезвоним в течении рабочего дня, чтобы обговорить детали проекта.")
$rng1.Font.HighlightColorIndex = "darkCyan"

# 3. Paragraph "Вы – истинный победитель!..." - highlight whole paragraph
#    (incl. paragraph mark) darkCyan
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*истинный победитель*") {
        $p.Range.Font.HighlightColorIndex = "darkCyan"
    }
}
